# "global totals-seperate user function"
# Rename the three team labels, update their "Current" totals, and move
# the active selection to C8 (matching the author's last interaction).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename teams (shared string text updates)
$ws.Range("A2").Value = "Team Pat"
$ws.Range("A3").Value = "Team Themba"
$ws.Range("A4").Value = "Team Sbu"

# Update "Current" totals for each team
$ws.Range("C2").Value = 5
$ws.Range("C3").Value = 10
$ws.Range("C4").Value = 13

# Move/save the active selection as it was when the workbook was saved
$ws.Range("C8").Select()
